$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = 0, bold font, thin border on all sides, centered horizontally, top vertically
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B1").VerticalAlignment = -4160    # xlTop
$ws.Range("B1").Borders.LineStyle = 1        # xlContinuous
$ws.Range("B1").Borders.Weight = 2           # xlThin

# A2 gets the same style as B1 - copy formatting to avoid a duplicate style entry
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("A2").Value = 0

# B2 = "disconnected_elements" (default/no special style)
$ws.Range("B2").Value = "disconnected_elements"
